$d = $word.ActiveDocument

# --------------------------------------------------------------------
# 0. The paragraph currently reads "ATM specificatio" + a "_GoBack"
#    bookmark + "n ". Relocate the (singleton) _GoBack bookmark out of
#    the way first so that fixing the typo below yields one clean run.
# --------------------------------------------------------------------
$docEnd = $d.Content.End
$tmpRange = $d.Range($docEnd - 1, $docEnd - 1)
$d.Bookmarks.Add("_GoBack", $tmpRange) | Out-Null

# --------------------------------------------------------------------
# 1. Fix the heading typo: "ATM specificatio" + "n " -> "ATM specification "
# --------------------------------------------------------------------
$d.Content.Find.Execute("ATM specificatio", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ATM specification", 2) | Out-Null

$heading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "ATM specification*") {
        $heading = $cand
        break
    }
}
$tailRange = $d.Range($heading.Range.Start + 17, $heading.Range.End - 1)
$tailRange.Text = " "

# Re-acquire the heading paragraph (now reads exactly "ATM specification ").
$heading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq "ATM specification `r") {
        $heading = $cand
        break
    }
}

# --------------------------------------------------------------------
# 2. Insert the two new introductory paragraphs right after the heading.
# --------------------------------------------------------------------
$insertPoint = $heading.Range
$insertPoint.Collapse(0)
$insertPoint.InsertParagraphAfter()
$insertPoint.InsertParagraphAfter()

$para2 = $heading.Next()
$para3 = $para2.Next()
$para2.Style = $d.Styles.Item("Normal")
$para3.Style = $d.Styles.Item("Normal")

# --------------------------------------------------------------------
# 3. Paragraph 2 text (single run, Times New Roman 13pt).
# --------------------------------------------------------------------
$r2 = $para2.Range
$r2.Collapse(1)
$r2.Font.Name = "Times New Roman"
$r2.Font.NameBi = "Times New Roman"
$r2.Font.Size = 13
$r2.Font.SizeBi = 13
$r2.InsertAfter("An ATM is primarily a machines and as we all know machines are made of component which differs according to the machines that we are talking about. ")

# --------------------------------------------------------------------
# 4. Paragraph 3 text: first half, then the relocated _GoBack bookmark,
#    then the second half ("be useless").
# --------------------------------------------------------------------
$r3 = $para3.Range
$r3.Collapse(1)
$r3.Font.Name = "Times New Roman"
$r3.Font.NameBi = "Times New Roman"
$r3.Font.Size = 13
$r3.Font.SizeBi = 13
$p3Start = $r3.Start
$firstHalf = "Out system" + [string][char]0x2019 + "s hardware is mainly a computer running some operating system (Windows OS) on which our software system will be installed. The computer without the software would "
$r3.InsertAfter($firstHalf)
$r3.Collapse(0)
$r3.InsertAfter("be useless")

$bmPos = $p3Start + $firstHalf.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
